$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_6_7_0"
$ws.Range("B2").Value = 0.2861104200961899
$ws.Range("C2").Value = -20.71836574409192
$ws.Range("D2").Value = 0.3524288140031163
$ws.Range("E2").Value = 0.06899947163376952
$ws.Range("F2").Value = 0.7900654673576355
$ws.Range("G2").Value = 0.832537055015564
$ws.Range("H2").Value = 1.340180039405823
$ws.Range("I2").Value = 1.071427702903748

$ws.Range("A3").Value = "model_6_7_14"
$ws.Range("B3").Value = 0.2936355747539209
$ws.Range("C3").Value = -28.63330599719765
$ws.Range("D3").Value = 0.4169816092555553
$ws.Range("E3").Value = -0.01594678367434144
$ws.Range("F3").Value = 0.7817373275756836
$ws.Range("G3").Value = 1.135943055152893
$ws.Range("H3").Value = 1.206584811210632
$ws.Range("I3").Value = 1.16918671131134

$ws.Range("A4").Value = "model_6_7_16"
$ws.Range("B4").Value = 0.293646072380482
$ws.Range("C4").Value = -28.64976759569252
$ws.Range("D4").Value = 0.4172702995976494
$ws.Range("E4").Value = -0.01599230460750256
$ws.Range("F4").Value = 0.7817255854606628
$ws.Range("G4").Value = 1.136574149131775
$ws.Range("H4").Value = 1.205987453460693
$ws.Range("I4").Value = 1.169239044189453

$ws.Range("A5").Value = "model_6_7_12"
$ws.Range("B5").Value = 0.2936573808639946
$ws.Range("C5").Value = -28.63894932779877
$ws.Range("D5").Value = 0.4172922906227431
$ws.Range("E5").Value = -0.01578255374672821
$ws.Range("F5").Value = 0.7817131280899048
$ws.Range("G5").Value = 1.136159420013428
$ws.Range("H5").Value = 1.205941915512085
$ws.Range("I5").Value = 1.168997764587402

$ws.Range("A6").Value = "model_6_7_8"
$ws.Range("B6").Value = 0.2936672991305826
$ws.Range("C6").Value = -28.61881462435281
$ws.Range("D6").Value = 0.4170458351364441
$ws.Range("E6").Value = -0.01563569913822316
$ws.Range("F6").Value = 0.7817021608352661
$ws.Range("G6").Value = 1.135387659072876
$ws.Range("H6").Value = 1.206451892852783
$ws.Range("I6").Value = 1.168828845024109

$ws.Range("A7").Value = "model_6_7_13"
$ws.Range("B7").Value = 0.2936696952189817
$ws.Range("C7").Value = -28.63319621447905
$ws.Range("D7").Value = 0.4171845355096531
$ws.Range("E7").Value = -0.01577270809494791
$ws.Range("F7").Value = 0.781699538230896
$ws.Range("G7").Value = 1.135938882827759
$ws.Range("H7").Value = 1.206164836883545
$ws.Range("I7").Value = 1.168986320495605

$ws.Range("A8").Value = "model_6_7_6"
$ws.Range("B8").Value = 0.293674237002948
$ws.Range("C8").Value = -28.61239865751223
$ws.Range("D8").Value = 0.4169567252692404
$ws.Range("E8").Value = -0.0155988093147108
$ws.Range("F8").Value = 0.7816944718360901
$ws.Range("G8").Value = 1.135141611099243
$ws.Range("H8").Value = 1.206636428833008
$ws.Range("I8").Value = 1.16878616809845

$ws.Range("A9").Value = "model_6_7_15"
$ws.Range("B9").Value = 0.2936783231671208
$ws.Range("C9").Value = -28.63045433520266
$ws.Range("D9").Value = 0.417020857645982
$ws.Range("E9").Value = -0.01586221608339855
$ws.Range("F9").Value = 0.7816900014877319
$ws.Range("G9").Value = 1.135833740234375
$ws.Range("H9").Value = 1.206503629684448
$ws.Range("I9").Value = 1.169089317321777

$ws.Range("A10").Value = "model_6_7_7"
$ws.Range("B10").Value = 0.2936819070626263
$ws.Range("C10").Value = -28.61186371646156
$ws.Range("D10").Value = 0.4169670539688231
$ws.Range("E10").Value = -0.01558053167844897
$ws.Range("F10").Value = 0.781686007976532
$ws.Range("G10").Value = 1.13512122631073
$ws.Range("H10").Value = 1.206614971160889
$ws.Range("I10").Value = 1.168765068054199

$ws.Range("A11").Value = "model_6_7_10"
$ws.Range("B11").Value = 0.293685701863993
$ws.Range("C11").Value = -28.61565631148485
$ws.Range("D11").Value = 0.4170242099855269
$ws.Range("E11").Value = -0.01559823202190813
$ws.Range("F11").Value = 0.7816817164421082
$ws.Range("G11").Value = 1.135266423225403
$ws.Range("H11").Value = 1.206496715545654
$ws.Range("I11").Value = 1.168785691261292

$ws.Range("A12").Value = "model_6_7_9"
$ws.Range("B12").Value = 0.2936879327453993
$ws.Range("C12").Value = -28.61199398288403
$ws.Range("D12").Value = 0.4169572762743158
$ws.Range("E12").Value = -0.01558955257097128
$ws.Range("F12").Value = 0.7816793322563171
$ws.Range("G12").Value = 1.135126233100891
$ws.Range("H12").Value = 1.206635236740112
$ws.Range("I12").Value = 1.16877555847168

$ws.Range("A13").Value = "model_6_7_11"
$ws.Range("B13").Value = 0.2937432894296195
$ws.Range("C13").Value = -28.60943777429184
$ws.Range("D13").Value = 0.4170695280031137
$ws.Range("E13").Value = -0.01545136136509528
$ws.Range("F13").Value = 0.781618058681488
$ws.Range("G13").Value = 1.135028123855591
$ws.Range("H13").Value = 1.206402778625488
$ws.Range("I13").Value = 1.168616652488708

$ws.Range("A14").Value = "model_6_7_17"
$ws.Range("B14").Value = 0.2937557870186885
$ws.Range("C14").Value = -28.61402749529649
$ws.Range("D14").Value = 0.4169045450822727
$ws.Range("E14").Value = -0.01567287579028065
$ws.Range("F14").Value = 0.7816042304039001
$ws.Range("G14").Value = 1.135203957557678
$ws.Range("H14").Value = 1.206744313240051
$ws.Range("I14").Value = 1.168871521949768

$ws.Range("A15").Value = "model_6_7_24"
$ws.Range("B15").Value = 0.2937639237727891
$ws.Range("C15").Value = -28.61143324553998
$ws.Range("D15").Value = 0.4168648533592783
$ws.Range("E15").Value = -0.01565979146857499
$ws.Range("F15").Value = 0.7815952301025391
$ws.Range("G15").Value = 1.135104656219482
$ws.Range("H15").Value = 1.206826448440552
$ws.Range("I15").Value = 1.168856620788574

$ws.Range("A16").Value = "model_6_7_22"
$ws.Range("B16").Value = 0.2937639237727891
$ws.Range("C16").Value = -28.61143324553998
$ws.Range("D16").Value = 0.4168648533592783
$ws.Range("E16").Value = -0.01565979146857499
$ws.Range("F16").Value = 0.7815952301025391
$ws.Range("G16").Value = 1.135104656219482
$ws.Range("H16").Value = 1.206826448440552
$ws.Range("I16").Value = 1.168856620788574

$ws.Range("A17").Value = "model_6_7_23"
$ws.Range("B17").Value = 0.2937639237727891
$ws.Range("C17").Value = -28.61143324553998
$ws.Range("D17").Value = 0.4168648533592783
$ws.Range("E17").Value = -0.01565979146857499
$ws.Range("F17").Value = 0.7815952301025391
$ws.Range("G17").Value = 1.135104656219482
$ws.Range("H17").Value = 1.206826448440552
$ws.Range("I17").Value = 1.168856620788574

$ws.Range("A18").Value = "model_6_7_21"
$ws.Range("B18").Value = 0.2937639237727891
$ws.Range("C18").Value = -28.61143324553998
$ws.Range("D18").Value = 0.4168648533592783
$ws.Range("E18").Value = -0.01565979146857499
$ws.Range("F18").Value = 0.7815952301025391
$ws.Range("G18").Value = 1.135104656219482
$ws.Range("H18").Value = 1.206826448440552
$ws.Range("I18").Value = 1.168856620788574

$ws.Range("A19").Value = "model_6_7_20"
$ws.Range("B19").Value = 0.2937639237727891
$ws.Range("C19").Value = -28.61143324553998
$ws.Range("D19").Value = 0.4168648533592783
$ws.Range("E19").Value = -0.01565979146857499
$ws.Range("F19").Value = 0.7815952301025391
$ws.Range("G19").Value = 1.135104656219482
$ws.Range("H19").Value = 1.206826448440552
$ws.Range("I19").Value = 1.168856620788574

$ws.Range("A20").Value = "model_6_7_18"
$ws.Range("B20").Value = 0.2937639237727891
$ws.Range("C20").Value = -28.61143324553998
$ws.Range("D20").Value = 0.4168648533592783
$ws.Range("E20").Value = -0.01565979146857499
$ws.Range("F20").Value = 0.7815952301025391
$ws.Range("G20").Value = 1.135104656219482
$ws.Range("H20").Value = 1.206826448440552
$ws.Range("I20").Value = 1.168856620788574

$ws.Range("A21").Value = "model_6_7_19"
$ws.Range("B21").Value = 0.2937639237727891
$ws.Range("C21").Value = -28.61143324553998
$ws.Range("D21").Value = 0.4168648533592783
$ws.Range("E21").Value = -0.01565979146857499
$ws.Range("F21").Value = 0.7815952301025391
$ws.Range("G21").Value = 1.135104656219482
$ws.Range("H21").Value = 1.206826448440552
$ws.Range("I21").Value = 1.168856620788574

$ws.Range("A22").Value = "model_6_7_4"
$ws.Range("B22").Value = 0.2938021818403691
$ws.Range("C22").Value = -28.56644883987209
$ws.Range("D22").Value = 0.4168013558586805
$ws.Range("E22").Value = -0.01492105553264511
$ws.Range("F22").Value = 0.7815529108047485
$ws.Range("G22").Value = 1.133380174636841
$ws.Range("H22").Value = 1.206957936286926
$ws.Range("I22").Value = 1.168006181716919

$ws.Range("A23").Value = "model_6_7_5"
$ws.Range("B23").Value = 0.293959056626457
$ws.Range("C23").Value = -28.52486020470613
$ws.Range("D23").Value = 0.4163024135802945
$ws.Range("E23").Value = -0.01460780359429514
$ws.Range("F23").Value = 0.7813792824745178
$ws.Range("G23").Value = 1.131785869598389
$ws.Range("H23").Value = 1.207990407943726
$ws.Range("I23").Value = 1.167645812034607

$ws.Range("A24").Value = "model_6_7_3"
$ws.Range("B24").Value = 0.2940249804557391
$ws.Range("C24").Value = -28.50114154415257
$ws.Range("D24").Value = 0.4168263442478815
$ws.Range("E24").Value = -0.01374759917611157
$ws.Range("F24").Value = 0.7813063263893127
$ws.Range("G24").Value = 1.130876779556274
$ws.Range("H24").Value = 1.206906199455261
$ws.Range("I24").Value = 1.166655778884888

$ws.Range("A25").Value = "model_6_7_2"
$ws.Range("B25").Value = 0.2947497588672933
$ws.Range("C25").Value = -28.25215298862607
$ws.Range("D25").Value = 0.4148217339790006
$ws.Range("E25").Value = -0.01105386850985091
$ws.Range("F25").Value = 0.7805042266845703
$ws.Range("G25").Value = 1.121332287788391
$ws.Range("H25").Value = 1.211054801940918
$ws.Range("I25").Value = 1.163555860519409

$ws.Range("A26").Value = "model_6_7_1"
$ws.Range("B26").Value = 0.2953004926839563
$ws.Range("C26").Value = -28.03010746835568
$ws.Range("D26").Value = 0.4135456326766912
$ws.Range("E26").Value = -0.008216852243094985
$ws.Range("F26").Value = 0.7798947095870972
$ws.Range("G26").Value = 1.112820386886597
$ws.Range("H26").Value = 1.213695764541626
$ws.Range("I26").Value = 1.160290837287903
